# Add a new "fm35" test row (calcrule 19) to the ftests sheet, following the
# same pattern as the existing rows, and select the newly added cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# Row 40 inherits the formatting of row 39 (its immediate predecessor).
$ws.Range("B39:I39").Copy()
$ws.Range("B40:I40").PasteSpecial(-4122)

# Column E (Allocrule) on this new row was typed as a plain, unformatted
# number - unlike row 39's right-aligned E cell - so drop its style.
$ws.Range("E40").ClearFormats()

$ws.Range("B40").Value = "fm35"
$ws.Range("C40").Value = "% Loss deductible with min and max deductible. Calcrule 19"
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 19
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = "in progress"
$ws.Range("I40").Value = "in progress"

# Row 41 got a stray formatted-but-empty B cell below the new entry.
$ws.Range("B40").Copy()
$ws.Range("B41").PasteSpecial(-4122)
$ws.Range("B41").ClearContents()

$excel.CutCopyMode = $false

[void]$ws.Range("H40:I40").Select()
